# Auto-generated edit script: updates the cryptos list (Price / Volume(1h) columns,
# plus a handful of row re-orderings where two coins swapped rank) to match the
# "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is unambiguously text (coin names, URLs, "##.##.##"-style
# prices with more than one dot, and the "  +x.xx%  " volume strings) -- safe to
# assign directly since Excel will not reinterpret them as a number.
$textCells = @(
    "D2|57.997.80",
    "E2|  +3.79%  ",
    "D3|2.457.78",
    "E3|  +3.09%  ",
    "E4|  -0.17%  ",
    "E5|  +8.14%  ",
    "E6|  +4.36%  ",
    "E7|  -0.83%  ",
    "E8|  +23.41%  ",
    "D9|2.483.24",
    "E9|  +3.95%  ",
    "E10|  +13.96%  ",
    "E11|  +6.06%  ",
    "E12|  +3.96%  ",
    "E13|  +1.79%  ",
    "D14|2.869.16",
    "E14|  +1.89%  ",
    "D15|57.865.72",
    "E15|  +2.82%  ",
    "E16|  +4.45%  ",
    "E17|  +3.06%  ",
    "D18|2.458.13",
    "E18|  +2.09%  ",
    "E19|  +6.87%  ",
    "E20|  +5.02%  ",
    "E21|  +4.78%  ",
    "B22|Uniswap",
    "C22|https://coinranking.com/coin/_H5FVG9iW+uniswap-uni",
    "E22|  +5.66%  ",
    "B23|Dai",
    "C23|https://coinranking.com/coin/MoTuySvg7+dai-dai",
    "E23|  -0.03%  ",
    "E24|  +3.00%  ",
    "E25|  +3.71%  ",
    "E26|  -1.26%  ",
    "E27|  +2.92%  ",
    "D28|2.538.02",
    "E28|  +0.25%  ",
    "E29|  +2.61%  ",
    "D30|0.0₃0809",
    "E30|  +5.54%  ",
    "E31|  -0.39%  ",
    "E32|  +6.30%  ",
    "E33|  +1.66%  ",
    "E34|  +4.92%  ",
    "E35|  +8.46%  ",
    "E36|  +5.28%  ",
    "E37|  +6.74%  ",
    "E38|  -0.72%  ",
    "E39|  +5.22%  ",
    "B40|Filecoin",
    "C40|https://coinranking.com/coin/ymQub4fuB+filecoin-fil",
    "E40|  +6.41%  ",
    "B41|OKB",
    "C41|https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb",
    "E41|  +2.72%  ",
    "E42|  +7.28%  ",
    "E43|  +10.62%  ",
    "E44|  -0.94%  ",
    "E45|  +3.76%  ",
    "E46|  +0.38%  ",
    "E47|  +5.04%  ",
    "B48|RenderToken",
    "C48|https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr",
    "E48|  +3.77%  ",
    "B49|WhiteBITCoin",
    "C49|https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt",
    "E49|  +0.28%  ",
    "E50|  +7.17%  ",
    "E51|  +11.87%  "
)

# Cells whose new value LOOKS like a plain decimal number (e.g. "1.00", "5.97").
# Excel/COM would silently coerce a bare assignment to a numeric cell (and mangle
# trailing zeros, e.g. "1.00" -> 1). The source workbook stores these as literal
# text, so force the cell to Text format first, then assign the string.
$numericLookingCells = @(
    "D4|1.00",
    "D5|159.33",
    "D6|494.70",
    "D7|0.993",
    "D8|0.612",
    "D10|6.30",
    "D12|0.336",
    "D16|21.15",
    "D19|4.74",
    "D20|327.25",
    "D21|10.18",
    "D22|5.97",
    "D23|0.997",
    "D24|58.62",
    "D25|0.408",
    "D26|0.993",
    "D27|0.162",
    "D29|7.43",
    "D31|0.997",
    "D32|18.93",
    "D33|151.10",
    "D35|5.41",
    "D37|3.82",
    "D38|0.839",
    "D39|1.41",
    "D40|3.59",
    "D41|34.40",
    "D43|281.05",
    "D44|0.991",
    "D45|0.608",
    "D46|0.0544",
    "D47|0.0232",
    "D48|4.74",
    "D49|10.24",
    "D50|18.22",
    "D51|0.691"
)

foreach ($item in $textCells) {
    $parts = $item.Split("|")
    $ref = $parts[0]
    $val = $parts[1]
    $ws.Range($ref).Value = $val
}

foreach ($item in $numericLookingCells) {
    $parts = $item.Split("|")
    $ref = $parts[0]
    $val = $parts[1]
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

